# Adds a new "2022-Q4" sheet (fund holdings) right after "总计", and a new
# summary row on "总计" for the 2022-Q4 quarter. All the other quarter
# sheets (2021-Q4 / 2021-Q2 / 2021-Q1 / 2020-Q4) are left untouched, just
# shifted one position to the right.

$wb = $excel.ActiveWorkbook

# Helper: assign a value as TEXT (keeps numeric-looking strings like "001917"
# or "5.91" from being coerced to a number) while leaving the cell's style
# untouched (matches the un-styled data cells used throughout these sheets).
function Set-TextCell($ws, $addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" worksheet by duplicating "2021-Q4" (so it
#    inherits the exact same layout / header styling), placed right after
#    "总计".
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$ws2021Q4 = $wb.Worksheets.Item("2021-Q4")
$ws2021Q4.Copy($null, $wsTotal)
$wsNew = $wb.Worksheets.Item("2021-Q4 (2)")
$wsNew.Name = "2022-Q4"

# Overwrite its two data rows with the 2022-Q4 fund holdings.
Set-TextCell $wsNew "B2" "001917"
Set-TextCell $wsNew "C2" "招商量化精选股票A"
Set-TextCell $wsNew "D2" "5.91"
Set-TextCell $wsNew "E2" "94.08"
Set-TextCell $wsNew "F2" "1.39"
Set-TextCell $wsNew "G2" "0.0821"
$wsNew.Range("H2").Value = 9

Set-TextCell $wsNew "B3" "007950"
Set-TextCell $wsNew "C3" "招商量化精选股票C"
Set-TextCell $wsNew "D3" "5.28"
Set-TextCell $wsNew "E3" "94.08"
Set-TextCell $wsNew "F3" "1.39"
Set-TextCell $wsNew "G3" "0.0734"
$wsNew.Range("H3").Value = 9

# ---------------------------------------------------------------------
# 2. Insert a new row 2 on "总计" for the 2022-Q4 summary entry, pushing
#    the existing rows (2021-Q4 / 2021-Q2 / 2021-Q1 / 2020-Q4) down by one.
# ---------------------------------------------------------------------
$wsTotal.Rows(2).Insert()

# Restore the row-2 formatting to match the other data rows (row 3, which
# now holds the original row-2 formatting/content after the insert shifted
# it down) instead of the header-derived formatting Insert() applied.
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$wsTotal.Range("B2:D2").Style = "Normal"

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.16

# Renumber the shifted rows' index column.
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("A5").Value = 3
$wsTotal.Range("A6").Value = 4
